# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2404   (left / "before" block of the AHB diff)
#   *_new -> *_FV2410   (right / "after" block of the AHB diff)
# Then expose the data range as a native Excel Table ("Table1") and
# freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells (row 1) ------------------------------------
# Columns A-J described the "before" (FV2404) input file, columns L-U the
# "after" (FV2410) input file; column K ("diff") is unaffected.
$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $fv2404Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2404Headers[$i]
}

for ($i = 0; $i -lt $fv2410Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2410Headers[$i]
}

# --- 2. Turn the used range into a native Excel table -----------------------
$dataRange = $ws.Range("A1:U89")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
# The source workbook does not bind a named table style (tableStyleInfo has
# no "name" attribute) - clear the default style Excel assigns on creation.
$tbl.TableStyle = ""

# --- 3. Freeze the header row ------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
